# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" detail table (rows 16-26, columns B:G) is re-sorted in
# ascending order by the "Periodo Mora" column (E). Each row's values
# (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico) travel
# together as a unit; only their row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "1909", 23333, 3500000),
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "1910", 140000, 3500000),
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "1911", 140000, 3500000),
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "1912", 140000, 3500000),
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "2001", 140000, 3500000),
    @("CC", "1047439864", "HERNANDO JOSE GARCIA DIAZ", "2002", 140000, 3500000),
    @("CC", "1143372563", "OSCAR DE JESUS SUAREZ TORRES", "2003", 50000, 2500000),
    @("CC", "1050970190", "LEIDER JAVIER TORRES RUIZ", "2007", 2341, 877803),
    @("CC", "1050970190", "LEIDER JAVIER TORRES RUIZ", "2008", 35112, 877803),
    @("CC", "1140886812", "ALEJANDRO DANIEL MANTILLA PAEZ", "2203", 3333, 2500000),
    @("CC", "1143358314", "CRISTIHAN CAMILO SUAREZ CANAVERAL", "2205", 2667, 1000000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}
